$wb = $excel.ActiveWorkbook

# --- Sheet "Beth S.": reorder In Stock list (A2:C5) by vendor (Grocery, Bath, Pet) ---
$wsBeth = $wb.Worksheets.Item("Beth S.")

$wsBeth.Cells.Item(2,1).Value = "Fruits"
$wsBeth.Cells.Item(2,2).Value = "Grocery"
$wsBeth.Cells.Item(2,3).Value = 7

$wsBeth.Cells.Item(3,1).Value = "Vegetables"
$wsBeth.Cells.Item(3,2).Value = "Grocery"
$wsBeth.Cells.Item(3,3).Value = 5

$wsBeth.Cells.Item(4,1).Value = "Body Butter"
$wsBeth.Cells.Item(4,2).Value = "Bath"
$wsBeth.Cells.Item(4,3).Value = 15

$wsBeth.Cells.Item(5,1).Value = "Premium Cat Food"
$wsBeth.Cells.Item(5,2).Value = "Pet"
$wsBeth.Cells.Item(5,3).Value = 11.99

# Visited list (G2:G4) also reordered to match the store visit order
$wsBeth.Cells.Item(2,7).Value = "Grocery"
$wsBeth.Cells.Item(3,7).Value = "Bath"
$wsBeth.Cells.Item(4,7).Value = "Pet"

# --- Sheet "Chris K.": reorder In Stock list (A2:C5) by vendor (Bath, Pet, Pet, Pet) ---
$wsChris = $wb.Worksheets.Item("Chris K.")

$wsChris.Cells.Item(2,1).Value = "Oatmeal Soap"
$wsChris.Cells.Item(2,2).Value = "Bath"
$wsChris.Cells.Item(2,3).Value = 7

$wsChris.Cells.Item(3,1).Value = "Cat Litter"
$wsChris.Cells.Item(3,2).Value = "Pet"
$wsChris.Cells.Item(3,3).Value = 9.99

$wsChris.Cells.Item(4,1).Value = "Premium Cat Food"
$wsChris.Cells.Item(4,2).Value = "Pet"
$wsChris.Cells.Item(4,3).Value = 11.99

$wsChris.Cells.Item(5,1).Value = "Brush"
$wsChris.Cells.Item(5,2).Value = "Pet"
$wsChris.Cells.Item(5,3).Value = 4.99

# Visited list (G2:G3) also reordered
$wsChris.Cells.Item(2,7).Value = "Bath"
$wsChris.Cells.Item(3,7).Value = "Pet"

# --- Switch the active/selected tab from "Mary M." to "Chris K." ---
$wsChris.Activate()
